$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 314.26666
$ws.Range("I5").Value = 276.25
$ws.Range("J5").Value = 466.33334
$ws.Range("K5").Value = 276.25
$ws.Range("L5").Value = 466.33334
$ws.Range("M5").Value = -161.25
$ws.Range("N5").Value = -696.33334
$ws.Range("H40").Value = 3846.4546
$ws.Range("J40").Value = 1363.3334
$ws.Range("L40").Value = 1363.3334
$ws.Range("N40").Value = -1713.3334
$ws.Range("H43").Value = 15394.909
$ws.Range("J43").Value = 6921
$ws.Range("L43").Value = 6921
$ws.Range("N43").Value = -7059

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5497.3687
$ws.Range("I32").Value = 3305.4365
$ws.Range("K32").Value = 3305.4365
$ws.Range("M32").Value = -3018.4365
$ws.Range("H45").Value = 11289.8
$ws.Range("I45").Value = 13549.875
$ws.Range("K45").Value = 13549.875
$ws.Range("M45").Value = -13172.875
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H132").Value = 3240.0857
$ws.Range("I132").Value = 3103.121
$ws.Range("K132").Value = 9309.363000000001
$ws.Range("M132").Value = -6779.363000000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 27771
$ws.Range("J27").Value = 27771
$ws.Range("L27").Value = 27771
$ws.Range("N27").Value = -28155
$ws.Range("H94").Value = 2000.5
$ws.Range("I94").Value = 1034
$ws.Range("J94").Value = 4900
$ws.Range("K94").Value = 1034
$ws.Range("L94").Value = 4900
$ws.Range("M94").Value = -583
$ws.Range("N94").Value = -5802
$ws.Range("H134").Value = 1808.7595
$ws.Range("I134").Value = 1806.3077
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5418.9231
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -2883.9231
$ws.Range("N134").Value = -11070

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 80827.08
$ws.Range("I31").Value = 114553.664
$ws.Range("J31").Value = 4942.25
$ws.Range("K31").Value = 114553.664
$ws.Range("L31").Value = 4942.25
$ws.Range("M31").Value = -114258.664
$ws.Range("N31").Value = -5532.25
$ws.Range("H34").Value = 80827.08
$ws.Range("I34").Value = 114553.664
$ws.Range("J34").Value = 4942.25
$ws.Range("K34").Value = 114553.664
$ws.Range("L34").Value = 4942.25
$ws.Range("M34").Value = -114351.664
$ws.Range("N34").Value = -5346.25
$ws.Range("H58").Value = 2524.3877
$ws.Range("I58").Value = 2195.525
$ws.Range("K58").Value = 2195.525
$ws.Range("M58").Value = -1992.525
$ws.Range("H80").Value = 33394
$ws.Range("J80").Value = 33394
$ws.Range("L80").Value = 33394
$ws.Range("N80").Value = -35640
$ws.Range("H83").Value = 33394
$ws.Range("J83").Value = 33394
$ws.Range("L83").Value = 100182
$ws.Range("N83").Value = -111414
$ws.Range("H107").Value = 1444.5555
$ws.Range("I107").Value = 1421.5714
$ws.Range("J107").Value = 1459.1818
$ws.Range("K107").Value = 1421.5714
$ws.Range("L107").Value = 1459.1818
$ws.Range("M107").Value = 498.4286
$ws.Range("N107").Value = -5299.1818
$ws.Range("H134").Value = 11570.889
$ws.Range("I134").Value = 6631.3076
$ws.Range("K134").Value = 19893.9228
$ws.Range("M134").Value = -17358.9228
$ws.Range("H136").Value = 2524.3877
$ws.Range("I136").Value = 2195.525
$ws.Range("K136").Value = 6586.575000000001
$ws.Range("M136").Value = -4036.575000000001
$ws.Range("H141").Value = 365755.16
$ws.Range("J141").Value = 365755.16
$ws.Range("L141").Value = 365755.16
$ws.Range("N141").Value = -376115.16

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 764
$ws.Range("I18").Value = 696.75
$ws.Range("J18").Value = 1033
$ws.Range("K18").Value = 2090.25
$ws.Range("L18").Value = 3099
$ws.Range("M18").Value = -1921.25
$ws.Range("N18").Value = -3437

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2848.611
$ws.Range("J80").Value = 2649.75
$ws.Range("L80").Value = 2649.75
$ws.Range("N80").Value = -4645.75
$ws.Range("H83").Value = 2848.611
$ws.Range("J83").Value = 2649.75
$ws.Range("L83").Value = 13248.75
$ws.Range("N83").Value = -23232.75
$ws.Range("H102").Value = 111112300
$ws.Range("I102").Value = 956.7143
$ws.Range("K102").Value = 956.7143
$ws.Range("M102").Value = 665.2857
$ws.Range("H111").Value = 99697
$ws.Range("J111").Value = 99697
$ws.Range("L111").Value = 99697
$ws.Range("N111").Value = -105831
$ws.Range("H113").Value = 2398.2917
$ws.Range("I113").Value = 2668.1765
$ws.Range("J113").Value = 1742.8572
$ws.Range("K113").Value = 2668.1765
$ws.Range("L113").Value = 1742.8572
$ws.Range("M113").Value = -498.1765
$ws.Range("N113").Value = -6082.8572

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4381
$ws.Range("I22").Value = 1350
$ws.Range("J22").Value = 5593.4
$ws.Range("K22").Value = 1350
$ws.Range("L22").Value = 5593.4
$ws.Range("M22").Value = -1055
$ws.Range("N22").Value = -6183.4
$ws.Range("H27").Value = 4381
$ws.Range("I27").Value = 1350
$ws.Range("J27").Value = 5593.4
$ws.Range("K27").Value = 1350
$ws.Range("L27").Value = 5593.4
$ws.Range("M27").Value = -1243
$ws.Range("N27").Value = -5807.4
$ws.Range("H68").Value = 296141.12
$ws.Range("J68").Value = 2001079.6
$ws.Range("L68").Value = 2001079.6
$ws.Range("N68").Value = -2002577.6
$ws.Range("H71").Value = 296141.12
$ws.Range("J71").Value = 2001079.6
$ws.Range("L71").Value = 10005398
$ws.Range("N71").Value = -10012886
$ws.Range("H93").Value = 1655.7778
$ws.Range("I93").Value = 1655.7778
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1655.7778
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -407.7778000000001
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 8285.571
$ws.Range("I122").Value = 7750
$ws.Range("J122").Value = 8999.666999999999
$ws.Range("K122").Value = 23250
$ws.Range("L122").Value = 26999.001
$ws.Range("M122").Value = -20800
$ws.Range("N122").Value = -31899.001
$ws.Range("H132").Value = 3030.2258
$ws.Range("I132").Value = 2346.1924
$ws.Range("K132").Value = 7038.5772
$ws.Range("M132").Value = -4508.5772
$ws.Range("H136").Value = 3379.9048
$ws.Range("I136").Value = 2804.3333
$ws.Range("K136").Value = 8412.999899999999
$ws.Range("M136").Value = -5862.999899999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1436.2307
$ws.Range("I96").Value = 1083.1428
$ws.Range("K96").Value = 1083.1428
$ws.Range("M96").Value = 289.8571999999999
$ws.Range("H132").Value = 1626.3966
$ws.Range("I132").Value = 1400.0212
$ws.Range("K132").Value = 4200.063599999999
$ws.Range("M132").Value = -1670.063599999999

Write-Host "All updates applied."
